$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.153.15"
$ws.Range("D3").Value = "1.904.70"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'306.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.5229"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.69%  "
$ws.Range("D8").Value = "'0.3765"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").Value = "'0.07254"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").Value = "'21.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "'0.9055"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "'0.08501"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.21%  "
$ws.Range("D13").Value = "1.915.98"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").Value = "'96.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("D15").Value = "'5.292"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "'0.000008663"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D20").Value = "27.191.29"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "'5.087"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "2.152.19"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").Value = "'6.436"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'2.341"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("D26").Value = "'146.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").Value = "'18.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").Value = "'1.747"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").Value = "'115.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").Value = "'4.921"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'4.821"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "'0.09299"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'0.8023"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("D34").Value = "'0.05061"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("D35").Value = "'1.245"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("D36").Value = "'3.446"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.82%  "
$ws.Range("D37").Value = "'2.949"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").Value = "'2.604"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").Value = "'0.5716"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("D40").Value = "'0.02001"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").Value = "'1.076"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'9.124"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").Value = "'6.637"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "'116.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").Value = "'0.1516"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").Value = "'0.4858"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("D47").Value = "'10.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").Value = "'0.9995"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "'1.622"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "'64.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.19%  "
